# The document carried two copies of the Pearson Edexcel footer logo
# (PearsonLogo.png) that were both mislabeled "image2.png" - rename them
# to "image1.png". The BTec header logo (BTec_Logo-Orange) was mislabeled
# "image1.jpg" - rename it to "image2.jpg". Pure name/label fix-up; the
# embedded picture bytes, sizing, and everything else stay untouched.
#
# InlineShapes inside headers/footers don't reliably accept a direct
# ".Name =" write in this object model (the handle can't be re-addressed
# for an inline anchor in those stories), so each picture is briefly
# converted to a floating shape - where the rename is applied - and then
# converted straight back to an inline shape, which restores the
# original <wp:inline> layout.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeByAltText($range, $altText, $newName) {
    for ($i = 1; $i -le $range.InlineShapes.Count; $i++) {
        $ishp = $range.InlineShapes.Item($i)
        if ($ishp.AlternativeText -eq $altText) {
            $floating = $ishp.ConvertToShape()
            $floating.Name = $newName
            $floating.ConvertToInlineShape() | Out-Null
        }
    }
}

# -- Footers: the Pearson Edexcel PowerPoint logo, image2.png -> image1.png
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        Rename-InlineShapeByAltText $ftr.Range "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image1.png"
    }
}

# -- Headers: the BTec_Logo-Orange picture, image1.jpg -> image2.jpg
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        Rename-InlineShapeByAltText $hdr.Range "BTec_Logo-Orange" "image2.jpg"
    }
}
